# Auto-generated PowerShell Excel COM-interop script
# Updates numeric result values in Sheet1 (pl_mw.xlsx / Case_5_224 line results)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"=0.5123330574875524; "C"=0.1191298221591452; "D"=0.04629170208878008; "E"=0.1009989519193013; "F"=0.9232754284974547; "I"=0.8623522309721174; "K"=0.3285408004581996; "L"=0.2043011126019394; "O"=3.289596083056196 }
  3 = @{ "B"=0.4690635668912364; "C"=0.1177010559591807; "D"=0.04408252612759611; "E"=0.1005822379082062; "F"=0.9252923456370326; "I"=0.8696967313908743; "K"=0.2894826279242011; "L"=0.1968817526097695; "O"=3.310028172141813 }
  4 = @{ "B"=0.4425834031724776; "C"=0.1168192702539699; "D"=0.04271260761978368; "E"=0.1003781762352993; "F"=0.9270721256038499; "I"=0.8746426740282551; "K"=0.2654701821498975; "L"=0.1924331771718357; "O"=3.324370843595332 }
  5 = @{ "B"=0.4318152270989515; "C"=0.1164588272108986; "D"=0.04215099352227725; "E"=0.1003080642353993; "F"=0.9279335626858796; "I"=0.876767926181941; "K"=0.2556778000435287; "L"=0.190647320518039; "O"=3.330667545682701 }
  6 = @{ "B"=0.430028571929455; "C"=0.1163989098053122; "D"=0.04205753573005211; "E"=0.1002972106128119; "F"=0.928084828309558; "I"=0.8771274509293043; "K"=0.2540513695668665; "L"=0.1903524119767752; "O"=3.331740404910732 }
  7 = @{ "B"=0.4424380869089646; "C"=0.1168144136428069; "D"=0.04270504706907019; "E"=0.1003771778375473; "F"=0.9270831918951785; "I"=0.8746708915794237; "K"=0.2653381467110876; "L"=0.1924089831580744; "O"=3.324453933202605 }
  8 = @{ "B"=0.4973959915839998; "C"=0.1186381373422165; "D"=0.04553279084149864; "E"=0.1008445281495263; "F"=0.9238585112105682; "I"=0.86479404598035; "K"=0.3150802177987089; "L"=0.2017207683886255; "O"=3.296268040646154 }
  9 = @{ "B"=0.6058367801742008; "C"=0.1221775609763398; "D"=0.05097010280109515; "E"=0.1021713803384721; "F"=0.9218304246864335; "I"=0.8488880017013827; "K"=0.4123613429557906; "L"=0.2208272049867048; "O"=3.255257870112246 }
  10 = @{ "B"=0.6858895072660971; "C"=0.1247543254560171; "D"=0.05489818443249561; "E"=0.1033957758698918; "F"=0.9229603830097588; "I"=0.8393126654137504; "K"=0.4836531176305527; "L"=0.2353788483306829; "O"=3.233827271308058 }
  11 = @{ "B"=0.722385307784208; "C"=0.1259211991945008; "D"=0.05667050661836726; "E"=0.1040068576514059; "F"=0.9240437183541133; "I"=0.8354149840221297; "K"=0.5160424854632311; "L"=0.2421101631597224; "O"=3.225968220368173 }
  12 = @{ "B"=0.7362161328598802; "C"=0.1263622761533014; "D"=0.05733951892219835; "E"=0.104246022671969; "F"=0.9245358177941512; "I"=0.8340049254962807; "K"=0.5283010296031989; "L"=0.2446751393164277; "O"=3.223264026484543 }
  13 = @{ "B"=0.7332369553322451; "C"=0.1262673180628369; "D"=0.0571955302693965; "E"=0.1041941692876946; "F"=0.9244261945103531; "I"=0.8343056751670161; "K"=0.5256612364970294; "L"=0.2441220168105218; "O"=3.223834330453514 }
  14 = @{ "B"=0.7235229688530467; "C"=0.125957502917224; "D"=0.056725589480358; "E"=0.1040263784487507; "F"=0.9240825630160359; "I"=0.8352976566218544; "K"=0.5170511392270498; "L"=0.2423208656733635; "O"=3.225740295690258 }
  15 = @{ "B"=0.7175742369856835; "C"=0.1257676280814763; "D"=0.05643745918342091; "E"=0.1039246119808652; "F"=0.923882739800213; "I"=0.8359138578320753; "K"=0.5117763272510558; "L"=0.2412196861464508; "O"=3.226943161381456 }
  16 = @{ "B"=0.6835059536810775; "C"=0.1246779581646393; "D"=0.05478206295090615; "E"=0.1033569273259012; "F"=0.9229010403202125; "I"=0.839576609847299; "K"=0.4815355056018973; "L"=0.2349411811621849; "O"=3.23437891586687 }
  17 = @{ "B"=0.6626259621274642; "C"=0.124008100175935; "D"=0.05376277428632648; "E"=0.1030225144908954; "F"=0.9224446124916739; "I"=0.8419409580603769; "K"=0.4629726385188064; "L"=0.2311180674697937; "O"=3.239424615939384 }
  18 = @{ "B"=0.6506238634292743; "C"=0.1236223175768103; "D"=0.05317513572823884; "E"=0.1028352621047048; "F"=0.922235673804785; "I"=0.8433439940478991; "K"=0.4522918881024509; "L"=0.2289296326022878; "O"=3.24250464347989 }
  19 = @{ "B"=0.6465614738901877; "C"=0.1234916135818835; "D"=0.05297593693371994; "E"=0.1027727369326428; "F"=0.922174134349838; "I"=0.8438264440853338; "K"=0.4486749220790784; "L"=0.2281904755096065; "O"=3.243578034833888 }
  20 = @{ "B"=0.6648479000901659; "C"=0.124079459457036; "D"=0.05387142138810219; "E"=0.1030575863259422; "F"=0.9224876540219; "I"=0.8416848056750617; "K"=0.4649490935112226; "L"=0.2315239564142075; "O"=3.238869082480164 }
  21 = @{ "B"=0.7263759189740995; "C"=0.1260485248846663; "D"=0.05686368047020096; "E"=0.1040754521786482; "F"=0.9241812741799009; "I"=0.8350044987010321; "K"=0.5195803177523715; "L"=0.2428494746225454; "O"=3.225173088543499 }
  22 = @{ "B"=0.7666498039224052; "C"=0.1273307884702035; "D"=0.05880687897975179; "E"=0.104785918712075; "F"=0.9257653160496915; "I"=0.8310227090928421; "K"=0.5552461858649167; "L"=0.250344418626895; "O"=3.217806608061721 }
  23 = @{ "B"=0.7451494730674142; "C"=0.1266468539490546; "D"=0.05777090356400549; "E"=0.1044025964409272; "F"=0.9248762237614585; "I"=0.833112705668583; "K"=0.5362143955426291; "L"=0.2463357401641701; "O"=3.22159320915739 }
  24 = @{ "B"=0.6638433553811751; "C"=0.1240472000192128; "D"=0.05382230712622516; "E"=0.1030417147550224; "F"=0.9224680283973328; "I"=0.8418004759198112; "K"=0.4640555652966327; "L"=0.2313404243502077; "O"=3.23911968100424 }
  25 = @{ "B"=0.5764318291403754; "C"=0.1212241234455718; "D"=0.04951081336775331; "E"=0.1017685615904149; "F"=0.9219190703958162; "I"=0.8528203267250802; "K"=0.3860744843158557; "L"=0.2155679991493287; "O"=3.264824820660493 }
}

foreach ($rowKey in $data.Keys) {
    $rowData = $data[$rowKey]
    foreach ($colKey in $rowData.Keys) {
        $cellRef = "$colKey$rowKey"
        $ws.Range($cellRef).Value = $rowData[$colKey]
    }
}
